$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Problem" value in row 2 (was "prevention instead of aftercare",
# now capitalized "Prevention instead of aftercare") so the table's
# filter/lookup matches correctly.
$ws.Range("F2").Value = "Prevention instead of aftercare"

# Restore the view: no frozen/scrolled topLeftCell, selection on A2.
$ws.Range("A2").Select()
